$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '66.068.85'
$ws.Range('E2').Value = '  -0.08%  '
Set-TextValue 'D3' '3.161.95'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue 'D5' '601.91'
$ws.Range('E5').Value = '  -0.42%  '
Set-TextValue 'D6' '153.96'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +2.47%  '
Set-TextValue 'D9' '3.159.23'
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('E10').Value = '  -2.07%  '
Set-TextValue 'D11' '5.52'
$ws.Range('E11').Value = '  -10.65%  '
Set-TextValue 'D12' '0.520'
$ws.Range('E12').Value = '  +2.05%  '
$ws.Range('E13').Value = '  -2.40%  '
Set-TextValue 'D14' '38.26'
$ws.Range('E14').Value = '  -0.95%  '
Set-TextValue 'D15' '3.679.78'
$ws.Range('E15').Value = '  -1.64%  '
Set-TextValue 'D16' '66.090.39'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('E17').Value = '  -1.00%  '
Set-TextValue 'D18' '3.161.69'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('E19').Value = '  +1.17%  '
Set-TextValue 'D20' '509.01'
$ws.Range('E20').Value = '  -0.57%  '
Set-TextValue 'D21' '15.37'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('E23').Value = '  -0.14%  '
Set-TextValue 'D24' '14.61'
$ws.Range('E24').Value = '  -4.43%  '
Set-TextValue 'D25' '84.43'
$ws.Range('E25').Value = '  -0.95%  '
Set-TextValue 'D26' '1.00'
$ws.Range('E26').Value = '  -0.06%  '
Set-TextValue 'D27' '2.98'
$ws.Range('E27').Value = '  -0.95%  '
Set-TextValue 'D28' '9.03'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('E29').Value = '  +5.50%  '
Set-TextValue 'D30' '3.05'
$ws.Range('E30').Value = '  +6.59%  '
Set-TextValue 'D31' '6.96'
$ws.Range('E31').Value = '  +2.35%  '
Set-TextValue 'D32' '27.92'
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  -2.20%  '
$ws.Range('E35').Value = '  -2.20%  '
Set-TextValue 'D36' '494.15'
$ws.Range('E36').Value = '  +2.74%  '
Set-TextValue 'D37' '54.62'
$ws.Range('E37').Value = '  -1.32%  '
Set-TextValue 'D38' '0.0881'
$ws.Range('E38').Value = '  -2.80%  '
Set-TextValue 'D39' '0.0419'
$ws.Range('E39').Value = '  -0.64%  '
Set-TextValue 'D40' '0.129'
$ws.Range('E40').Value = '  +7.68%  '
Set-TextValue 'D41' '8.75'
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('E42').Value = '  +4.70%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D43' '2.80'
$ws.Range('E43').Value = '  -5.96%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D44' '0.294'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('E45').Value = '  -1.90%  '
Set-TextValue 'D46' '2.816.38'
$ws.Range('E46').Value = '  -4.32%  '
Set-TextValue 'D47' '27.77'
$ws.Range('E47').Value = '  -3.92%  '
Set-TextValue 'D49' '2.36'
$ws.Range('E49').Value = '  +2.17%  '
Set-TextValue 'D50' '0.117'
$ws.Range('E50').Value = '  +0.64%  '
Set-TextValue 'D51' '34.89'
$ws.Range('E51').Value = '  +3.15%  '
